$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.059.17'
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = '1.829.63'
$ws.Range("E3").Value = '  -0.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '241.42'
$ws.Range("D6").Value = '0.6342'
$ws.Range("E6").Value = '  -4.63%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.80'
$ws.Range("E8").Value = '  +6.88%  '
$ws.Range("D9").Value = '0.2937'
$ws.Range("E9").Value = '  +0.81%  '
$ws.Range("D10").Value = '0.07344'
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("D11").Value = '22.88'
$ws.Range("E11").Value = '  +1.37%  '
$ws.Range("D12").Value = '0.07682'
$ws.Range("D13").Value = '1.829.74'
$ws.Range("E13").Value = '  -0.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.990'
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("D15").Value = '0.6635'
$ws.Range("E15").Value = '  -0.20%  '
$ws.Range("D16").Value = '81.98'
$ws.Range("E16").Value = '  -1.98%  '
$ws.Range("D17").Value = '6.068'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008662'
$ws.Range("E18").Value = '  +4.78%  '
$ws.Range("D19").Value = '28.982.57'
$ws.Range("E19").Value = '  -0.67%  '
$ws.Range("D20").Value = '2.079.81'
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("D21").Value = '12.42'
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").Value = '224.17'
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("D24").Value = '7.125'
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = '1.001'
$ws.Range("E25").Value = '  -0.13%  '
$ws.Range("E26").Value = '  -1.57%  '
$ws.Range("D27").Value = '8.466'
$ws.Range("E27").Value = '  -1.79%  '
$ws.Range("D28").Value = '0.1373'
$ws.Range("E28").Value = '  -1.27%  '
$ws.Range("D29").Value = '17.88'
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").Value = '1.501'
$ws.Range("E30").Value = '  -0.67%  '
$ws.Range("D31").Value = '4.096'
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("D32").Value = '4.026'
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("E33").Value = '  +0.81%  '
$ws.Range("D34").Value = '0.05299'
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").Value = '1.831'
$ws.Range("E35").Value = '  -1.68%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.7391'
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("D37").Value = '1.154'
$ws.Range("E37").Value = '  +2.48%  '
$ws.Range("D38").Value = '2.652'
$ws.Range("E38").Value = '  -1.20%  '
$ws.Range("D39").Value = '1.292.79'
$ws.Range("E39").Value = '  -1.48%  '
$ws.Range("E40").Value = '  +1.22%  '
$ws.Range("D41").Value = '0.01783'
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("D42").Value = '6.289'
$ws.Range("E42").Value = '  +5.90%  '
$ws.Range("D43").Value = '0.8957'
$ws.Range("E43").Value = '  -2.08%  '
$ws.Range("D44").Value = '0.9996'
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("D45").Value = '102.71'
$ws.Range("E45").Value = '  +0.91%  '
$ws.Range("D46").Value = '1.978.33'
$ws.Range("E46").Value = '  +0.40%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '0.5136'
$ws.Range("E47").Value = '  -0.54%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '64.09'
$ws.Range("E48").Value = '  +1.53%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000120'
$ws.Range("E49").Value = '  -7.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.730'
$ws.Range("E50").Value = '  -1.65%  '
$ws.Range("D51").Value = '0.05819'
$ws.Range("E51").Value = '  -1.85%  '
